# Apply cryptos price/volume updates per the commit diff.
# For column D values that are plain decimal numbers (e.g. "212.37"),
# prefix with a single-quote so Excel stores them as TEXT (matching the
# original inlineStr text cells) instead of auto-converting to a Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "88.467.18"
$ws.Range("E2").Value = "  +0.68%  "

$ws.Range("D3").Value = "3.279.22"
$ws.Range("E3").Value = "  -1.60%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'212.37"
$ws.Range("E5").Value = "  -2.86%  "

$ws.Range("D6").Value = "'628.20"
$ws.Range("E6").Value = "  -1.39%  "

$ws.Range("D7").Value = "'0.375"
$ws.Range("E7").Value = "  +14.33%  "

$ws.Range("D8").Value = "'0.721"
$ws.Range("E8").Value = "  +17.16%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").Value = "3.271.89"
$ws.Range("E10").Value = "  -2.00%  "

$ws.Range("E11").Value = "  -5.82%  "

$ws.Range("D12").Value = "'0.187"
$ws.Range("E12").Value = "  +12.03%  "

$ws.Range("E13").Value = "  -4.78%  "

$ws.Range("D14").Value = "'5.49"
$ws.Range("E14").Value = "  +1.35%  "

$ws.Range("D15").Value = "'34.17"
$ws.Range("E15").Value = "  -0.44%  "

$ws.Range("D16").Value = "3.882.54"
$ws.Range("E16").Value = "  -1.88%  "

$ws.Range("D17").Value = "88.293.96"
$ws.Range("E17").Value = "  +0.84%  "

$ws.Range("D18").Value = "3.304.15"
$ws.Range("E18").Value = "  -0.99%  "

$ws.Range("E19").Value = "  -0.97%  "

$ws.Range("D20").Value = "'14.05"
$ws.Range("E20").Value = "  -3.46%  "

$ws.Range("D21").Value = "'436.96"
$ws.Range("E21").Value = "  -2.63%  "

$ws.Range("D22").Value = "'8.91"
$ws.Range("E22").Value = "  -2.21%  "

$ws.Range("D23").Value = "'5.35"
$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("D24").Value = "'7.45"
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("B25").Value = "Aptos"
$ws.Range("C25").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D25").Value = "'12.29"
$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "'5.24"
$ws.Range("E26").Value = "  -2.43%  "

$ws.Range("E27").Value = "  -1.58%  "

$ws.Range("D28").Value = "'77.12"
$ws.Range("E28").Value = "  -2.10%  "

$ws.Range("E29").Value = "  +4.29%  "

$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("D31").Value = "'0.181"
$ws.Range("E31").Value = "  -3.59%  "

$ws.Range("D32").Value = "'0.997"
$ws.Range("E32").Value = "  -0.27%  "

$ws.Range("D33").Value = "'570.63"
$ws.Range("E33").Value = "  -5.23%  "

$ws.Range("D34").Value = "'8.77"
$ws.Range("E34").Value = "  -5.80%  "

$ws.Range("D35").Value = "'1.38"
$ws.Range("E35").Value = "  -10.78%  "

$ws.Range("B36").Value = "PancakeSwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D36").Value = "'1.96"
$ws.Range("E36").Value = "  -4.29%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "'7.13"
$ws.Range("E37").Value = "  +6.22%  "

$ws.Range("E38").Value = "  -7.52%  "

$ws.Range("D39").Value = "'22.79"
$ws.Range("E39").Value = "  -2.79%  "

$ws.Range("D40").Value = "'21.82"
$ws.Range("E40").Value = "  +2.08%  "

$ws.Range("E41").Value = "  -0.30%  "

$ws.Range("D42").Value = "'3.10"
$ws.Range("E42").Value = "  +1.26%  "

$ws.Range("D43").Value = "'0.402"
$ws.Range("E43").Value = "  -4.08%  "

$ws.Range("E44").Value = "  -1.56%  "

$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("D46").Value = "'153.74"
$ws.Range("E46").Value = "  -2.90%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.137"
$ws.Range("E47").Value = "  +21.07%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'180.86"
$ws.Range("E48").Value = "  -4.57%  "

$ws.Range("D49").Value = "'44.83"
$ws.Range("E49").Value = "  -2.82%  "

$ws.Range("E50").Value = "  -3.99%  "

$ws.Range("E51").Value = "  -1.26%  "
